$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete column A ("کد") entirely - header cell + shared string row shifts left
$ws.Columns.Item(1).Delete()

# 2. Move the "quick access" column (now at column S after the delete) to be column C
#    (right after "نام و نام خانوادگی"). Insert a blank column first, copy the quick-access
#    column's content/format into it, then remove the now-duplicated source column -- this
#    avoids the stray zero-width column records that Cut+Insert leaves behind.
$ws.Columns.Item(3).Insert()
$ws.Columns.Item(20).Copy() | Out-Null
$ws.Columns.Item(3).PasteSpecial(-4104) | Out-Null   # xlPasteAll
$ws.Columns.Item(20).Delete()

# 3. Update text of the (now relocated) quick-access header cell
$ws.Range("C1").Value = "دسترسی سریع (۱) فعال"

# 4. Apply uniform header style across A1:S1 - bold font, centered, yellow fill
$headerRange = $ws.Range("A1:S1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4108    # xlCenter
$headerRange.Interior.Color = 65535       # RGB(255,255,0) yellow

# reset reading order to context (remove forced RTL on A1) for all header cells
$headerRange.ReadingOrder = -5002         # xlReadingOrderContext

# 5. Set custom widths for the two columns that didn't exist before (columns A and B keep
#    the widths they inherited from the old B/C columns automatically).
$ws.Columns.Item(3).ColumnWidth = 18.75
$ws.Columns.Item(4).ColumnWidth = 13.75

# 6. Set active selection to B6
$ws.Range("B6").Select()
